# Edit the "headers" worksheet:
#  - insert a new row above the old "starshipit" row for a "starshipit_ui"
#    header block (cut-down / UI-facing report headers)
#  - append 3 new trailing headers to the "starshipit" row
#    (tracking_short_status, tracking_number, results.last_updated_date)
#  - give the category column (column A) a thin right border
#  - leave the selection where the author left it

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("headers")
$ws.Activate() | Out-Null

# --- insert the new "starshipit_ui" header row above row 3 ("starshipit") ---
$ws.Rows(3).Insert()

$ws.Range("A3").Value = "starshipit_ui"
$ws.Range("B3").Value = "Order Date"
$ws.Range("C3").Value = "Printed Date"
$ws.Range("D3").Value = "Delivered Date"
$ws.Range("E3").Value = "Price"
$ws.Range("F3").Value = "Carrier"
$ws.Range("G3").Value = "Package Sent"
$ws.Range("H3").Value = "AccountName"
$ws.Range("I3").Value = "Item Skus"
$ws.Range("B3:I3").Style = "Good"

# --- append the new trailing columns to the "starshipit" row (now row 4) ---
$ws.Range("BF4").Value = "tracking_short_status"
$ws.Range("BG4").Value = "tracking_number"
$ws.Range("BH4").Value = "results.last_updated_date"
$ws.Range("BF4:BH4").Style = "Good"

# --- add a thin right border under the category column (column A) ---
$catRange = $ws.Range("A1:A15")
$rightBorder = $catRange.Borders.Item(10)
$rightBorder.LineStyle = 1
$rightBorder.Weight = 2

# --- restore the author's final selection/view state ---
$ws.Range("BH26").Select() | Out-Null
